$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.146.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "'1.894.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.95%  "
$ws.Range("D4").Value = "'0.9975"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'247.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'0.9974"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.4989"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").Value = "'44.75"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").Value = "'0.2944"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.37%  "
$ws.Range("D10").Value = "'0.06656"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("D11").Value = "'1.890.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.97%  "
$ws.Range("D12").Value = "'17.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "'0.07196"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "'0.6786"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("D15").Value = "'85.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").Value = "'4.854"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").Value = "'30.132.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "'0.000008020"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.85%  "
$ws.Range("D19").Value = "'0.9964"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'12.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.22%  "
$ws.Range("D21").Value = "'2.135.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("D22").Value = "'0.9977"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'4.778"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.21%  "
$ws.Range("D24").Value = "'5.666"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.67%  "
$ws.Range("D25").Value = "'9.180"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.97%  "
$ws.Range("D26").Value = "'147.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").Value = "'133.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("D28").Value = "'16.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("D29").Value = "'1.950"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.74%  "
$ws.Range("D30").Value = "'1.371"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").Value = "'4.233"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").Value = "'0.08748"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").Value = "'3.954"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.05138"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.93%  "
$ws.Range("D35").Value = "'1.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").Value = "'0.7088"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.98%  "
$ws.Range("D37").Value = "'2.662"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "'2.776"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("D39").Value = "'2.235"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").Value = "'0.9406"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "'0.01664"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("D42").Value = "'6.076"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D43").Value = "'0.4228"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.14%  "
$ws.Range("D44").Value = "'0.9960"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "'103.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").Value = "'7.493"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.45%  "
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("D48").Value = "'0.05724"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("D49").Value = "'32.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.67%  "
$ws.Range("D50").Value = "'8.297"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("D51").Value = "'0.3747"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.47%  "
